$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column R (2020) to the table, mirroring the formatting of
# column Q (the previous year's column) for each row.
$ws.Range("Q4:Q14").Copy()
$ws.Range("R4:R14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("R4").Value = 2020
$ws.Range("R5").Value = 2.1
$ws.Range("R6").Value = 2.4
$ws.Range("R7").Value = 1.4
$ws.Range("R8").Value = 3.2
$ws.Range("R9").Value = 2.4
$ws.Range("R10").Value = 0.8
$ws.Range("R11").Value = 2.2000000000000002
$ws.Range("R12").Value = 4.5
$ws.Range("R13").Value = 1.4
$ws.Range("R14").Value = 3.2

# Update the active selection to match the saved workbook state
$ws.Range("R16:R17").Select()
